# DOMA-3100 add formatter convert to number for some colomns
#
# The ticket-analytics export template has per-ticket numeric columns
# (processing / completed / canceled / deferred / closed / new_or_reopened)
# for the two templated rows (row 2 -> tickets[i], row 3 -> tickets[i+1]).
# Those columns should render through the "formatN()" Carbone.io formatter
# and be stored with a numeric ("0") cell format instead of plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Switch the numeric data columns (C..H) on the two template rows to a
#    real number format ("0") instead of the inherited text format.
$ws.Range("C2:H3").NumberFormat = "0"

# 2) Append the ":formatN()" Carbone formatter to each of those template
#    placeholders so the exported values get converted to numbers.
$ws.Range("C2").Value = "{d.tickets[i].processing:formatN()}"
$ws.Range("D2").Value = "{d.tickets[i].completed:formatN()}"
$ws.Range("E2").Value = "{d.tickets[i].canceled:formatN()}"
$ws.Range("F2").Value = "{d.tickets[i].deferred:formatN()}"
$ws.Range("G2").Value = "{d.tickets[i].closed:formatN()}"
$ws.Range("H2").Value = "{d.tickets[i].new_or_reopened:formatN()}"

$ws.Range("C3").Value = "{d.tickets[i+1].processing:formatN()}"
$ws.Range("D3").Value = "{d.tickets[i+1].completed:formatN()}"
$ws.Range("E3").Value = "{d.tickets[i+1].canceled:formatN()}"
$ws.Range("F3").Value = "{d.tickets[i+1].deferred:formatN()}"
$ws.Range("G3").Value = "{d.tickets[i+1].closed:formatN()}"
$ws.Range("H3").Value = "{d.tickets[i+1].new_or_reopened:formatN()}"
